# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) values on Sheet1 were stored like "5-6-2013-14"
# (month-day-season) and need to become an actual game date "2014-05-06"
# (the NBA "2013-14" season game played on 5/6 actually occurred in May
# 2014), for rows 2 through 31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF$row")
    if ($cell.Value() -eq "5-6-2013-14") {
        # Assigning the literal text "2014-05-06" directly would cause
        # Excel to auto-recognize it as a date and silently convert the
        # cell to a date value/format. Instead, compute it as a string
        # formula result and then paste only the value back onto the
        # cell, which keeps it a plain text string without touching the
        # cell's number format/style.
        $cell.Formula = '="2014-05-06"'
        $cell.Copy()
        $cell.PasteSpecial(-4163)  # xlPasteValues
    }
}

$excel.CutCopyMode = 0
